$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("B2").Value = "2025-04-26 06:04:32"
$ws.Range("C2").Value = "John Smith found battery 3. Now John Smith is Tired"

# --- Add new row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "2025-04-26 06:40:03"
$c3 = @"
John Smith received battery 2 from Suppliers Battery New.
battery 2's state was New.
Thus John Smith carried out the following actions:
Store, 
Now John Smith is Frustrated.

"@
$ws.Range("C3").Value = $c3
$ws.Rows.Item(3).EntireRow.AutoFit()

# --- Add new row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "2025-04-26 06:50:20"
$c4 = @"
John Smith shipped battery 2 to Suppliers Never Death Row.
Now John Smith is Frustrated.

"@
$ws.Range("C4").Value = $c4
$ws.Rows.Item(4).EntireRow.AutoFit()

# --- Widen column C (closest reachable value to the target 211.2 chars) ---
$ws.Columns.Item(3).ColumnWidth = 210.3
